$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ROUTERS (existing sheet1) - trim to routers only, add a second router
# ---------------------------------------------------------------------
$routers = $wb.Worksheets.Item(1)
$routers.Name = "Routers"

# Row 4 held a switch ("Switch glowny4") that is moving to its own sheet.
$routers.Rows.Item(4).Delete()

# New router added as row 3.
$routers.Range("A3").Value = "Router V2"
$routers.Range("B3").Value = "Opis asdasd"
$routers.Range("C3").Value = "192.168.1.2"
$routers.Range("D3").Value = 3
$routers.Range("F3").Value = "Admin"
$routers.Range("G3").Value = "authPriv"
$routers.Range("H3").Value = "MD5"
$routers.Range("I3").Value = "Password"
$routers.Range("J3").Value = "DES"
$routers.Range("K3").Value = "Password"

# ---------------------------------------------------------------------
# Create both new sheets up front (positions/ids are fixed by THIS
# creation order: Services is created first -> sheetId 2, Switches is
# created second -> sheetId 3. Both get inserted right after Routers,
# so creating Switches last leaves the final tab order
# Routers, Switches, Services).
#
# NOTE: worksheet variables returned by Worksheets.Add()/Item() track a
# *position*, not an identity - once another sheet gets inserted in
# front of where they sit, the old variable silently starts resolving
# to whatever now occupies that slot. So every sheet is re-fetched by
# name (Worksheets.Item("Name")) right before it is touched, instead of
# reusing a variable captured before a later Add() could have shifted
# it.
# ---------------------------------------------------------------------
$wb.Worksheets.Add([System.Type]::Missing, $routers).Name = "Services"
$wb.Worksheets.Add([System.Type]::Missing, $routers).Name = "Switches"

# ---------------------------------------------------------------------
# SERVICES data
# ---------------------------------------------------------------------
$services = $wb.Worksheets.Item("Services")

$services.Range("A1").Value = "category"
$services.Range("B1").Value = "uptime"
$services.Range("C1").Value = "ping"
$services.Range("D1").Value = "interface_status"
$services.Range("E1").Value = "interface_utilization"
$services.Range("F1").Value = "chassis_temperature"
$services.Range("G1").Value = "fan_status"

$services.Range("A2").Value = "Routers"
$services.Range("B2").Value = $true
$services.Range("C2").Value = $true
$services.Range("D2").Value = "fa0/1"

$services.Range("A3").Value = "Switches"
$services.Range("B3").Value = $true
$services.Range("C3").Value = $true
$services.Range("D3").Value = "fa0/3"
$services.Range("E3").Value = $true
$services.Range("F3").Value = $true
$services.Range("G3").Value = $true

# Services carries an explicit print/page setup (portrait, paper size 9).
$services.PageSetup.PaperSize = 9
$services.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# SWITCHES data
# ---------------------------------------------------------------------
$switches = $wb.Worksheets.Item("Switches")

$switches.Range("A1").Value = "name"
$switches.Range("B1").Value = "description"
$switches.Range("C1").Value = "address"
$switches.Range("D1").Value = "snmp_version"
$switches.Range("E1").Value = "community"
$switches.Range("F1").Value = "security_name"
$switches.Range("G1").Value = "security_level"
$switches.Range("H1").Value = "auth_protocol"
$switches.Range("I1").Value = "priv_key"
$switches.Range("J1").Value = "priv_protocol"
$switches.Range("K1").Value = "auth_key"

$switches.Range("A2").Value = "Switch drugi"
$switches.Range("B2").Value = "Opis 2"
$switches.Range("C2").Value = "192.168.1.5"
$switches.Range("D2").Value = 3
$switches.Range("E2").Value = ""
$switches.Range("F2").Value = "Admin"
$switches.Range("G2").Value = "authPriv"
$switches.Range("H2").Value = "MD5"
$switches.Range("I2").Value = "Password"
$switches.Range("J2").Value = "DES"
$switches.Range("K2").Value = "ok"

$switches.Range("A3").Value = "Switch glowny4"
$switches.Range("B3").Value = "Opisdwa"
$switches.Range("C3").Value = "192.168.1.6"
$switches.Range("D3").Value = "2c"
$switches.Range("E3").Value = "Password"
$switches.Range("F3").Value = ""
$switches.Range("G3").Value = ""
$switches.Range("H3").Value = ""
$switches.Range("I3").Value = ""
$switches.Range("J3").Value = ""
$switches.Range("K3").Value = ""

# ---------------------------------------------------------------------
# View state: selections per sheet + which tab is active on open.
# Re-fetch every sheet by name here too, now that no further Add()
# calls will run and shift anybody's position.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Routers").Range("C7").Select()
$wb.Worksheets.Item("Switches").Range("D21").Select()
$wb.Worksheets.Item("Services").Range("E4").Select()
$wb.Worksheets.Item("Services").Activate()
